# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - VALOR MORA: 384000 -> 512000
# - Cant. Periodos: 3 -> 4
# - Novedad de Ingreso / Novedad de Retiro header columns swap
# - Worker/period table rebuilt: 2 workers (YESENIA / OSMIRO) x 4 periods
#   (2505, 2506, 2507, 2508) instead of x3 periods (2505, 2506, 2507)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the 2 extra data rows needed to go from 6 to 8 data rows ----
# (pushes the old last row 21 -> 23, and the footer block 26/27 -> 28/29)
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

# Give the 2 new rows the same formatting as the other plain data rows
# (row 20 is an ordinary, non-last data row).
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B20:J20").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Header values (text only changes, same wording) --------------------
$ws.Range("E11").Value = 512000
$ws.Range("F13").Value = 4

# Novedad de Ingreso / Novedad de Retiro swap positions in the header row
$ws.Range("H15").Value = "Novedad de Retiro"
$ws.Range("I15").Value = "Novedad de Ingreso"

# --- Rebuild the worker/period table (rows 16-23) ------------------------
# Tipo Doc (B) is "CC" for every row already; keep as-is, just make sure.
$rows = @(16,17,18,19,20,21,22,23)
foreach ($r in $rows) {
    $ws.Range("B$r").Value = "CC"
}

# Row 16: YESENIA / 2505
$ws.Range("C16").Value = "1047416130"
$ws.Range("D16").Value = "YESENIA ALZAMORA CASTAÑO"
$ws.Range("E16").Value = "2505"

# Row 17: OSMIRO / 2505
$ws.Range("C17").Value = "1128063438"
$ws.Range("D17").Value = "OSMIRO DE JESUS MORALES VASQUEZ"
$ws.Range("E17").Value = "2505"

# Row 18: YESENIA / 2506
$ws.Range("C18").Value = "1047416130"
$ws.Range("D18").Value = "YESENIA ALZAMORA CASTAÑO"
$ws.Range("E18").Value = "2506"

# Row 19: OSMIRO / 2506
$ws.Range("C19").Value = "1128063438"
$ws.Range("D19").Value = "OSMIRO DE JESUS MORALES VASQUEZ"
$ws.Range("E19").Value = "2506"

# Row 20: YESENIA / 2507
$ws.Range("C20").Value = "1047416130"
$ws.Range("D20").Value = "YESENIA ALZAMORA CASTAÑO"
$ws.Range("E20").Value = "2507"

# Row 21 (new): OSMIRO / 2507
$ws.Range("C21").Value = "1128063438"
$ws.Range("D21").Value = "OSMIRO DE JESUS MORALES VASQUEZ"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 64000
$ws.Range("G21").Value = 1600000

# Row 22 (new): YESENIA / 2508
$ws.Range("C22").Value = "1047416130"
$ws.Range("D22").Value = "YESENIA ALZAMORA CASTAÑO"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 64000
$ws.Range("G22").Value = 1600000

# Row 23 (old last row, now shifted down): OSMIRO / 2508
$ws.Range("C23").Value = "1128063438"
$ws.Range("D23").Value = "OSMIRO DE JESUS MORALES VASQUEZ"
$ws.Range("E23").Value = "2508"
